$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add header for new column P, matching the style of the other header cells (A1)
$ws.Range("P1").Value = "model_timestamp"
$ws.Range("A1").Copy()
$ws.Range("P1").PasteSpecial(-4122)

# Fill P2:P13 with the model_timestamp value 6993
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 16).Value = 6993
}
